$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: Excel auto-converts plain decimal-looking text into a Number,
# which would corrupt values like "11.50" (trailing zero) or "1.00" -> 1.
# Prefixing such values with a leading apostrophe forces Excel to keep them
# as literal text (matches how the source workbook stores them).
function Set-CellText {
    param($range, [string]$text)
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-CellText $ws.Range('D2') '95.515.62'
Set-CellText $ws.Range('E2') '  -0.88%  '
Set-CellText $ws.Range('D3') '3.463.38'
Set-CellText $ws.Range('E3') '  +3.95%  '
Set-CellText $ws.Range('E4') '  +0.14%  '
Set-CellText $ws.Range('D5') '241.07'
Set-CellText $ws.Range('E5') '  -3.41%  '
Set-CellText $ws.Range('D6') '644.41'
Set-CellText $ws.Range('E6') '  -1.41%  '
Set-CellText $ws.Range('D7') '1.48'
Set-CellText $ws.Range('E7') '  +6.81%  '
Set-CellText $ws.Range('D8') '0.403'
Set-CellText $ws.Range('E8') '  -4.08%  '
Set-CellText $ws.Range('E9') '  +0.11%  '
Set-CellText $ws.Range('D10') '0.998'
Set-CellText $ws.Range('E10') '  +0.45%  '
Set-CellText $ws.Range('D11') '3.462.39'
Set-CellText $ws.Range('E11') '  +3.98%  '
Set-CellText $ws.Range('E12') '  -3.64%  '
Set-CellText $ws.Range('D13') '41.71'
Set-CellText $ws.Range('E13') '  +3.59%  '
Set-CellText $ws.Range('D14') '6.13'
Set-CellText $ws.Range('E14') '  +0.85%  '
Set-CellText $ws.Range('D15') '95.269.27'
Set-CellText $ws.Range('E15') '  -0.84%  '
Set-CellText $ws.Range('D16') '4.114.03'
Set-CellText $ws.Range('E16') '  +4.19%  '
Set-CellText $ws.Range('E17') '  +2.87%  '
Set-CellText $ws.Range('D18') '8.46'
Set-CellText $ws.Range('E18') '  -0.29%  '
Set-CellText $ws.Range('D19') '3.479.87'
Set-CellText $ws.Range('E19') '  +4.51%  '
Set-CellText $ws.Range('D20') '17.94'
Set-CellText $ws.Range('E20') '  +5.27%  '
Set-CellText $ws.Range('D21') '11.50'
Set-CellText $ws.Range('E21') '  +9.30%  '
Set-CellText $ws.Range('D22') '0.512'
Set-CellText $ws.Range('E22') '  +0.46%  '
Set-CellText $ws.Range('D23') '503.86'
Set-CellText $ws.Range('E23') '  +0.11%  '
Set-CellText $ws.Range('D24') '3.18'
Set-CellText $ws.Range('E24') '  -4.77%  '
Set-CellText $ws.Range('D25') '0.0000192'
Set-CellText $ws.Range('E25') '  -2.48%  '
Set-CellText $ws.Range('D26') '6.61'
Set-CellText $ws.Range('E26') '  +1.13%  '
Set-CellText $ws.Range('D27') '95.10'
Set-CellText $ws.Range('E27') '  -0.90%  '
Set-CellText $ws.Range('B28') 'WrappedeETH'
Set-CellText $ws.Range('C28') 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-CellText $ws.Range('D28') '3.648.32'
Set-CellText $ws.Range('E28') '  +4.20%  '
Set-CellText $ws.Range('B29') 'Aptos'
Set-CellText $ws.Range('C29') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText $ws.Range('D29') '12.06'
Set-CellText $ws.Range('E29') '  +0.05%  '
Set-CellText $ws.Range('D30') '11.71'
Set-CellText $ws.Range('E30') '  +5.71%  '
Set-CellText $ws.Range('E31') '  +0.08%  '
Set-CellText $ws.Range('E32') '  +11.64%  '
Set-CellText $ws.Range('D33') '0.137'
Set-CellText $ws.Range('E33') '  -4.23%  '
Set-CellText $ws.Range('E34') '  -0.88%  '
Set-CellText $ws.Range('D35') '31.17'
Set-CellText $ws.Range('E35') '  +11.27%  '
Set-CellText $ws.Range('D36') '1.00'
Set-CellText $ws.Range('E36') '  -0.06%  '
Set-CellText $ws.Range('D37') '0.571'
Set-CellText $ws.Range('E37') '  +4.84%  '
Set-CellText $ws.Range('D38') '7.78'
Set-CellText $ws.Range('E38') '  +2.57%  '
Set-CellText $ws.Range('D39') '1.44'
Set-CellText $ws.Range('E39') '  -2.77%  '
Set-CellText $ws.Range('D40') '525.34'
Set-CellText $ws.Range('E40') '  +3.39%  '
Set-CellText $ws.Range('E41') '  +0.03%  '
Set-CellText $ws.Range('E42') '  +0.20%  '
Set-CellText $ws.Range('D43') '0.914'
Set-CellText $ws.Range('E43') '  +9.74%  '
Set-CellText $ws.Range('D44') '24.12'
Set-CellText $ws.Range('E44') '  -0.86%  '
Set-CellText $ws.Range('E45') '  +2.69%  '
Set-CellText $ws.Range('B46') 'Filecoin'
Set-CellText $ws.Range('C46') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws.Range('D46') '5.65'
Set-CellText $ws.Range('E46') '  +2.78%  '
Set-CellText $ws.Range('B47') 'VeChain'
Set-CellText $ws.Range('C47') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws.Range('D47') '0.0417'
Set-CellText $ws.Range('E47') '  +1.10%  '
Set-CellText $ws.Range('B48') 'Stacks'
Set-CellText $ws.Range('C48') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText $ws.Range('D48') '2.15'
Set-CellText $ws.Range('E48') '  +9.24%  '
Set-CellText $ws.Range('B49') 'MantraDAO'
Set-CellText $ws.Range('C49') 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-CellText $ws.Range('D49') '3.48'
Set-CellText $ws.Range('E49') '  -4.37%  '
Set-CellText $ws.Range('D50') '53.63'
Set-CellText $ws.Range('E50') '  +0.97%  '
Set-CellText $ws.Range('D51') '3.19'
Set-CellText $ws.Range('E51') '  +2.06%  '
